$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44477
$ws.Range("H2").Value = "Sin especificar"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 500
$ws.Range("K2").Value = 1400
$ws.Range("L2").Value = 1500
$ws.Range("M2").Value = 1460
$ws.Range("N2").Value = "$/kilo"
$ws.Range("O2").Value = "Provincia de Linares"
$ws.Range("P2").Value = 1460

$ws.Range("D3").Value = 44860
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 1100
$ws.Range("K3").Value = 1500
$ws.Range("L3").Value = 1700
$ws.Range("M3").Value = 1609
$ws.Range("N3").Value = "$/kilo"
$ws.Range("O3").Value = "Provincia de Linares"
$ws.Range("P3").Value = 1609

$ws.Range("D4").Value = 44496
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 550
$ws.Range("K4").Value = 1500
$ws.Range("L4").Value = 2000
$ws.Range("M4").Value = 1773
$ws.Range("N4").Value = "$/paquete"
$ws.Range("O4").Value = "Provincia de Linares"
$ws.Range("P4").Value = 1773

$ws.Range("D5").Value = 44868
$ws.Range("H5").Value = "Sin especificar"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 1000
$ws.Range("K5").Value = 1200
$ws.Range("L5").Value = 1300
$ws.Range("M5").Value = 1250
$ws.Range("N5").Value = "$/kilo"
$ws.Range("O5").Value = "Región del Maule"
$ws.Range("P5").Value = 1250

$ws.Range("D6").Value = 44868
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Segunda"
$ws.Range("J6").Value = 200
$ws.Range("K6").Value = 1000
$ws.Range("L6").Value = 1000
$ws.Range("M6").Value = 1000
$ws.Range("N6").Value = "$/kilo"
$ws.Range("O6").Value = "Región del Maule"
$ws.Range("P6").Value = 1000

$ws.Range("D7").Value = 44881
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 2600
$ws.Range("L7").Value = 2700
$ws.Range("M7").Value = 2650
$ws.Range("N7").Value = "$/kilo"
$ws.Range("O7").Value = "Provincia de Linares"
$ws.Range("P7").Value = 2650

$ws.Range("D8").Value = 44881
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Segunda"
$ws.Range("J8").Value = 100
$ws.Range("K8").Value = 2400
$ws.Range("L8").Value = 2400
$ws.Range("M8").Value = 2400
$ws.Range("N8").Value = "$/kilo"
$ws.Range("O8").Value = "Provincia de Linares"
$ws.Range("P8").Value = 2400

$ws.Range("D9").Value = 45245
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 300
$ws.Range("K9").Value = 1800
$ws.Range("L9").Value = 2000
$ws.Range("M9").Value = 1900
$ws.Range("N9").Value = "$/kilo"
$ws.Range("O9").Value = "Provincia de Linares"
$ws.Range("P9").Value = 1900

$ws.Range("D10").Value = 44875
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 300
$ws.Range("K10").Value = 1500
$ws.Range("L10").Value = 1600
$ws.Range("M10").Value = 1550
$ws.Range("N10").Value = "$/kilo"
$ws.Range("O10").Value = "Provincia de Linares"
$ws.Range("P10").Value = 1550

$ws.Range("D11").Value = 44519
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 250
$ws.Range("K11").Value = 1200
$ws.Range("L11").Value = 1300
$ws.Range("M11").Value = 1240
$ws.Range("N11").Value = "$/kilo"
$ws.Range("O11").Value = "Provincia de Linares"
$ws.Range("P11").Value = 1240

$ws.Range("D12").Value = 44489
$ws.Range("H12").Value = "Sin especificar"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 600
$ws.Range("K12").Value = 1400
$ws.Range("L12").Value = 1500
$ws.Range("M12").Value = 1450
$ws.Range("N12").Value = "$/kilo"
$ws.Range("O12").Value = "Provincia de Linares"
$ws.Range("P12").Value = 1450

$ws.Range("D13").Value = 45203
$ws.Range("H13").Value = "Sin especificar"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 400
$ws.Range("K13").Value = 1400
$ws.Range("L13").Value = 1500
$ws.Range("M13").Value = 1450
$ws.Range("N13").Value = "$/kilo"
$ws.Range("O13").Value = "Provincia de Linares"
$ws.Range("P13").Value = 1450

$ws.Range("D14").Value = 44468
$ws.Range("H14").Value = "Verde"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 500
$ws.Range("K14").Value = 1800
$ws.Range("L14").Value = 2000
$ws.Range("M14").Value = 1920
$ws.Range("N14").Value = "$/kilo"
$ws.Range("O14").Value = "Provincia de Linares"
$ws.Range("P14").Value = 1920

$ws.Range("D15").Value = 45246
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 300
$ws.Range("K15").Value = 2000
$ws.Range("L15").Value = 2000
$ws.Range("M15").Value = 2000
$ws.Range("N15").Value = "$/kilo"
$ws.Range("O15").Value = "Provincia de Linares"
$ws.Range("P15").Value = 2000

$ws.Range("D16").Value = 45246
$ws.Range("H16").Value = "Sin especificar"
$ws.Range("I16").Value = "Segunda"
$ws.Range("J16").Value = 300
$ws.Range("K16").Value = 1500
$ws.Range("L16").Value = 1500
$ws.Range("M16").Value = 1500
$ws.Range("N16").Value = "$/kilo"
$ws.Range("O16").Value = "Provincia de Linares"
$ws.Range("P16").Value = 1500

$ws.Range("D17").Value = 45230
$ws.Range("H17").Value = "Sin especificar"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 500
$ws.Range("K17").Value = 1500
$ws.Range("L17").Value = 1500
$ws.Range("M17").Value = 1500
$ws.Range("N17").Value = "$/kilo"
$ws.Range("O17").Value = "Provincia de Linares"
$ws.Range("P17").Value = 1500

$ws.Range("D18").Value = 44511
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 600
$ws.Range("K18").Value = 1300
$ws.Range("L18").Value = 1400
$ws.Range("M18").Value = 1350
$ws.Range("N18").Value = "$/kilo"
$ws.Range("O18").Value = "Provincia de Linares"
$ws.Range("P18").Value = 1350

$ws.Range("D19").Value = 44545
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 550
$ws.Range("K19").Value = 1700
$ws.Range("L19").Value = 1800
$ws.Range("M19").Value = 1755
$ws.Range("N19").Value = "$/kilo"
$ws.Range("O19").Value = "Provincia de Linares"
$ws.Range("P19").Value = 1755

$ws.Range("D20").Value = 44839
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 500
$ws.Range("K20").Value = 1700
$ws.Range("L20").Value = 1800
$ws.Range("M20").Value = 1760
$ws.Range("N20").Value = "$/kilo"
$ws.Range("O20").Value = "Provincia de Linares"
$ws.Range("P20").Value = 1760

$ws.Range("D21").Value = 44526
$ws.Range("H21").Value = "Sin especificar"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 100
$ws.Range("K21").Value = 1500
$ws.Range("L21").Value = 1600
$ws.Range("M21").Value = 1550
$ws.Range("N21").Value = "$/kilo"
$ws.Range("O21").Value = "Provincia de Linares"
$ws.Range("P21").Value = 1550

$ws.Range("D22").Value = 45202
$ws.Range("H22").Value = "Verde"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 1600
$ws.Range("L22").Value = 1600
$ws.Range("M22").Value = 1600
$ws.Range("N22").Value = "$/kilo"
$ws.Range("O22").Value = "Provincia de Linares"
$ws.Range("P22").Value = 1600

$ws.Range("D23").Value = 44510
$ws.Range("H23").Value = "Sin especificar"
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 600
$ws.Range("K23").Value = 1300
$ws.Range("L23").Value = 1400
$ws.Range("M23").Value = 1350
$ws.Range("N23").Value = "$/kilo"
$ws.Range("O23").Value = "Provincia de Linares"
$ws.Range("P23").Value = 1350

$ws.Range("D24").Value = 44876
$ws.Range("H24").Value = "Sin especificar"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 350
$ws.Range("K24").Value = 1500
$ws.Range("L24").Value = 1600
$ws.Range("M24").Value = 1557
$ws.Range("N24").Value = "$/kilo"
$ws.Range("O24").Value = "Provincia de Linares"
$ws.Range("P24").Value = 1557

$ws.Range("D25").Value = 44524
$ws.Range("H25").Value = "Sin especificar"
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 200
$ws.Range("K25").Value = 1500
$ws.Range("L25").Value = 1600
$ws.Range("M25").Value = 1550
$ws.Range("N25").Value = "$/kilo"
$ws.Range("O25").Value = "Provincia de Talca"
$ws.Range("P25").Value = 1550
